$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers M1, N1 (copy header style from existing header cell)
$ws.Range("L1").Copy() | Out-Null
$ws.Range("M1:N1").PasteSpecial(-4122) | Out-Null
$ws.Range("M1").Value = "Project Title"
$ws.Range("N1").Value = "Project Description"

# Update Country column D2:D6 from Delhi to India
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 4).Value = "India"
}

# Fix K4 tech skills ordering/content
$ws.Range("K4").Value = "Manufacturing,Manufacturing Operation Tools,Mechanical,Mechanical Processes"

# Populate Project Title (M) and Project Description (N) columns
$ws.Range("M2").Value = "Production management software for Army Base Workshop"
$ws.Range("N2").Value = @"
Production management software for Army Base Workshop
Project Description
Description: 
A. Existing setup : The OH of a tank is a 6 stage process spread over 144 days. It involves multiple entities working in tandem to achieve the target in time.
B. Problem Statement :
1. Development of a software application/ package to monitor the progress of production process.
2. The software should be able to identify crucial bottlenecks and suggest corrective action in advance
3. assist in easy HR management
Skill sets : Proficiency in programming languages, data structures, algorithms, problem solving and communication
"@

$ws.Range("M3").Value = "Engine Health Monitoring System."
$ws.Range("N3").Value = @"
Engine Health Monitoring System.
Project Description
Existing Setup : Instantaneous and not continuous monitoring, leads to reactive maintenance than preventive maintenance.
Problem Statement.
(i) Analog System to be digitized.
(ii) Requirement of ECU that takes input from all sensors.
(iii) Deviation of critical engine parameters from normal to be logged and highlighted.
(iv) Pin-point fault diagnosis.
(v) Predict and alert for preventive maintenance in advance.
(vi) Remote Diagnosis of Faults.
Skill Sets: Electro Mechanics
"@

$ws.Range("M4").Value = "Title:Re-manufacturing of certain internal components for reduction in weight thereby improving power to weight ratio."
$ws.Range("N4").Value = @"
Title:Re-manufacturing of certain internal components for reduction in weight thereby improving power to weight ratio.
Project Description
Description:
(a) Existing Setup: Metallic parts designed for ruggedness.
Therefore, they are bulky and reduce the power to weight ratio.
Additionally, the fuel/ammunition within is susceptible to catching fire
and even exploding due to enemy fire.
(b) Problem Statement: Finding suitable alternative material for .
re-manufacturing these components such that the new component is:-
(i) Light Weight-leading to increase in power to weight ratio.
(ii) Self Sealing in cases where fuel/oils/lubricants are contained
within.
(iii) Blast proof stowage for ammunition.
Skill Sets: Mechanical
"@

$ws.Range("M5").Value = "Re-designing of external components for reduction in thermal signature, thereby enhancing Tank protection."
$ws.Range("N5").Value = @"
Re-designing of external components for reduction in thermal signature, thereby enhancing Tank protection.
Project Description
Description:
(a) Existing Setup. Entirely metallic, designed ruggedness.
(b) Problem Statement. Finding suitable alternative
technology for reducing the external heat signature of the Tank.
Skill Sets:MechanicaI, Machine Design
"@

$ws.Range("M6").Value = "lndigenization of Induction Sensor of PWDM and Potentiometer (HYDR1 )"
$ws.Range("N6").Value = @"
lndigenization of Induction Sensor of PWDM and Potentiometer (HYDR1 )
Project Description
Description:
(a) Existing Setup: Both are essential for Turret Traverse in Stab
and non Stab mode and not available through Ordnance channel.
(b) Problem Statement:-
(i) Not available in normal source of supply.
(ii) UNSV in over hauling of Tanks.
(iii) Required 01(potentiometer (HYDR1)) and 02 (Induction
Sensor) qty per Tank.
(iv) Not developed by BEL.
(v) Ex-import spare.
(vi) Local Vendor/Supplier is not available.
Skill Sets:- Electro Mechanics
"@

for ($r = 1; $r -le 6; $r++) {
    $ws.Rows.Item($r).AutoFit() | Out-Null
}

Write-Host "done"
